$d = $word.ActiveDocument

$replacements = @(
    @{Old="303÷8="; New="335÷3="},
    @{Old="449÷3="; New="569÷6="},
    @{Old="875÷8="; New="673÷3="},
    @{Old="332÷9="; New="853÷5="},
    @{Old="322÷6="; New="277÷9="},
    @{Old="844÷8="; New="953÷7="},
    @{Old="946÷3="; New="326÷2="},
    @{Old="832÷4="; New="563÷6="},
    @{Old="950÷5="; New="122÷3="},
    @{Old="571÷3="; New="441÷6="},
    @{Old="947÷9="; New="626÷8="},
    @{Old="783÷8="; New="918÷4="},
    @{Old="191÷6="; New="305÷8="},
    @{Old="684÷2="; New="971÷5="},
    @{Old="903÷2="; New="452÷2="},
    @{Old="329÷2="; New="775÷5="},
    @{Old="831÷4="; New="116÷3="},
    @{Old="216÷7="; New="308÷8="},
    @{Old="364÷4="; New="314÷7="},
    @{Old="461÷9="; New="896÷5="},
    @{Old="324÷7="; New="173÷6="},
    @{Old="830÷2="; New="119÷3="},
    @{Old="367÷6="; New="379÷4="},
    @{Old="338÷6="; New="460÷4="},
    @{Old="438÷8="; New="421÷5="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
